$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update team-specific time-matrix probabilities (rows 2-19)

# Row 2
$ws.Range("B2").Value = 0.1944444444444444
$ws.Range("C2").Value = 0.5740740740740741
$ws.Range("J2").Value = 0.01234567901234568
$ws.Range("P2").Value = 0.1512345679012346
$ws.Range("S2").Value = 0.06790123456790123

# Row 3
$ws.Range("B3").Value = 0.005
$ws.Range("C3").Value = 0.02
$ws.Range("J3").Value = 0.05
$ws.Range("P3").Value = 0.6899999999999999
$ws.Range("S3").Value = 0.235

# Row 4
$ws.Range("P4").Value = 0.631578947368421
$ws.Range("S4").Value = 0.3684210526315789

# Row 6
$ws.Range("B6").Value = 0.05508474576271186
$ws.Range("D6").Value = 0.02542372881355932
$ws.Range("F6").Value = 0.0635593220338983
$ws.Range("J6").Value = 0.288135593220339
$ws.Range("O6").Value = 0.02542372881355932
$ws.Range("Q6").Value = 0.1483050847457627
$ws.Range("R6").Value = 0.05932203389830509
$ws.Range("S6").Value = 0.3347457627118644

# Row 7
$ws.Range("B7").Value = 0.09345794392523364
$ws.Range("D7").Value = 0.004672897196261682
$ws.Range("F7").Value = 0.0514018691588785
$ws.Range("J7").Value = 0.1588785046728972
$ws.Range("O7").Value = 0.01869158878504673
$ws.Range("Q7").Value = 0.1448598130841121
$ws.Range("R7").Value = 0.1168224299065421
$ws.Range("S7").Value = 0.411214953271028

# Row 8
$ws.Range("B8").Value = 0.1090487238979118
$ws.Range("D8").Value = 0.009280742459396751
$ws.Range("F8").Value = 0.07192575406032482
$ws.Range("J8").Value = 0.111368909512761
$ws.Range("O8").Value = 0.03248259860788863
$ws.Range("Q8").Value = 0.1670533642691415
$ws.Range("R8").Value = 0.1183294663573086
$ws.Range("S8").Value = 0.3805104408352668

# Row 9
$ws.Range("B9").Value = 0.125
$ws.Range("D9").Value = 0.0108695652173913
$ws.Range("F9").Value = 0.07608695652173914
$ws.Range("J9").Value = 0.1304347826086956
$ws.Range("O9").Value = 0.01630434782608696
$ws.Range("Q9").Value = 0.1684782608695652
$ws.Range("R9").Value = 0.07065217391304347
$ws.Range("S9").Value = 0.4021739130434783

# Row 10
$ws.Range("B10").Value = 0.1119293078055965
$ws.Range("D10").Value = 0.01840942562592047
$ws.Range("E10").Value = 0.001472754050073638
$ws.Range("F10").Value = 0.0625920471281296
$ws.Range("J10").Value = 0.1170839469808542
$ws.Range("O10").Value = 0.02135493372606775
$ws.Range("Q10").Value = 0.227540500736377
$ws.Range("R10").Value = 0.08100147275405008
$ws.Range("S10").Value = 0.3586156111929308

# Row 11
$ws.Range("G11").Value = 0.1394658753709199
$ws.Range("J11").Value = 0.1127596439169139
$ws.Range("K11").Value = 0.2017804154302671
$ws.Range("L11").Value = 0.5400593471810089
$ws.Range("S11").Value = 0.005934718100890208

# Row 12
$ws.Range("G12").Value = 0.7554347826086957
$ws.Range("J12").Value = 0.2065217391304348
$ws.Range("L12").Value = 0.0108695652173913
$ws.Range("S12").Value = 0.02717391304347826

# Row 13
$ws.Range("G13").Value = 0.6153846153846154
$ws.Range("J13").Value = 0.3076923076923077
$ws.Range("S13").Value = 0.07692307692307693

# Row 15
$ws.Range("F15").Value = 0.02262443438914027
$ws.Range("H15").Value = 0.1719457013574661
$ws.Range("I15").Value = 0.07692307692307693
$ws.Range("J15").Value = 0.2805429864253394
$ws.Range("K15").Value = 0.07239819004524888
$ws.Range("M15").Value = 0.03167420814479638
$ws.Range("O15").Value = 0.05429864253393665
$ws.Range("S15").Value = 0.2895927601809955

# Row 16
$ws.Range("F16").Value = 0.04455445544554455
$ws.Range("H16").Value = 0.1138613861386139
$ws.Range("I16").Value = 0.0594059405940594
$ws.Range("J16").Value = 0.4405940594059406
$ws.Range("K16").Value = 0.1138613861386139
$ws.Range("M16").Value = 0.0297029702970297
$ws.Range("O16").Value = 0.04455445544554455
$ws.Range("S16").Value = 0.1534653465346535

# Row 17
$ws.Range("F17").Value = 0.02536997885835095
$ws.Range("H17").Value = 0.1627906976744186
$ws.Range("I17").Value = 0.08668076109936575
$ws.Range("J17").Value = 0.4355179704016913
$ws.Range("K17").Value = 0.09725158562367865
$ws.Range("M17").Value = 0.0105708245243129
$ws.Range("O17").Value = 0.04228329809725159
$ws.Range("S17").Value = 0.1395348837209302

# Row 18
$ws.Range("F18").Value = 0.0187793427230047
$ws.Range("H18").Value = 0.1596244131455399
$ws.Range("I18").Value = 0.08450704225352113
$ws.Range("J18").Value = 0.4178403755868544
$ws.Range("K18").Value = 0.08450704225352113
$ws.Range("M18").Value = 0.0187793427230047
$ws.Range("O18").Value = 0.06103286384976526
$ws.Range("S18").Value = 0.1549295774647887

# Row 19
$ws.Range("F19").Value = 0.01864280387770321
$ws.Range("H19").Value = 0.1968680089485459
$ws.Range("I19").Value = 0.07233407904548844
$ws.Range("J19").Value = 0.3683818046234154
$ws.Range("K19").Value = 0.1200596569724087
$ws.Range("M19").Value = 0.0238627889634601
$ws.Range("N19").Value = 0.001491424310216256
$ws.Range("O19").Value = 0.06338553318419091
$ws.Range("S19").Value = 0.1349739000745712
